$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 6.1312172268804073
$ws.Range("C2").Value = 4.0099020443588804
$ws.Range("D2").Value = 4.9288288682325225
$ws.Range("E2").Value = 6.1312172268804064

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 5.6429209075745135
$ws.Range("C3").Value = 4.4647050986894543
$ws.Range("D3").Value = 4.4647050986894481
$ws.Range("E3").Value = 5.1644469555492307

# Update the selected range to match the new active selection in the sheet view
$ws.Range("B1:E3").Select()

$wb.Save()
